# Update "想去人数" (want-to-go count) values in column F across sheets
# per the scraped-data refresh (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 550
$ws1.Range("F5").Value = 6297
$ws1.Range("F6").Value = 705
$ws1.Range("F7").Value = 1079
$ws1.Range("F13").Value = 1130
$ws1.Range("F18").Value = 1405
$ws1.Range("F20").Value = 361
$ws1.Range("F21").Value = 384
$ws1.Range("F23").Value = 1059
$ws1.Range("F25").Value = 2165
$ws1.Range("F30").Value = 3483

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 143
$ws2.Range("F18").Value = 372
$ws2.Range("F32").Value = 1560

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 550
$ws4.Range("F12").Value = 143
$ws4.Range("F15").Value = 6297
$ws4.Range("F18").Value = 1079
$ws4.Range("F27").Value = 1130
$ws4.Range("F30").Value = 372
$ws4.Range("F34").Value = 1405
$ws4.Range("F36").Value = 361
$ws4.Range("F37").Value = 384
$ws4.Range("F42").Value = 2165
$ws4.Range("F46").Value = 1560
$ws4.Range("F50").Value = 3483
